$wb = $excel.ActiveWorkbook

# "studies" sheet - add PMID column in H1
$wsStudies = $wb.Worksheets.Item("studies")
$wsStudies.Range("H1").Value = "PMID"

# "counts" sheet - add notes column in F1
$wsCounts = $wb.Worksheets.Item("counts")
$wsCounts.Range("F1").Value = "notes"

# Update selections / active sheet
$wsStudies.Range("H2").Select()
$wsCounts.Range("F2").Select()
$wsCounts.Activate()
